$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.452.23'
$ws.Range("E2").Value = '  +2.11%  '

$ws.Range("D3").Value = '2.328.26'
$ws.Range("E3").Value = '  +0.16%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.50'
$ws.Range("E5").Value = '  +6.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.94'
$ws.Range("E6").Value = '  +2.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.992'
$ws.Range("E7").Value = '  -0.73%  '

$ws.Range("E8").Value = '  +0.65%  '

$ws.Range("D9").Value = '2.359.59'
$ws.Range("E9").Value = '  +1.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.103'
$ws.Range("E10").Value = '  +2.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.43'
$ws.Range("E11").Value = '  +3.20%  '

$ws.Range("E12").Value = '  +1.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +5.27%  '

$ws.Range("D14").Value = '2.753.49'
$ws.Range("E14").Value = '  +0.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.54'
$ws.Range("E15").Value = '  +0.61%  '

$ws.Range("D16").Value = '57.541.73'
$ws.Range("E16").Value = '  +2.29%  '

$ws.Range("E17").Value = '  +1.45%  '

$ws.Range("D18").Value = '2.384.24'
$ws.Range("E18").Value = '  +2.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.44'
$ws.Range("E19").Value = '  +5.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.57'
$ws.Range("E20").Value = '  +2.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.23'
$ws.Range("E21").Value = '  +2.31%  '

$ws.Range("E22").Value = '  +4.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.996'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.74'
$ws.Range("E24").Value = '  +0.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.170'
$ws.Range("E25").Value = '  +5.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.52'
$ws.Range("E26").Value = '  -0.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("E27").Value = '  -0.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.43'
$ws.Range("E28").Value = '  +10.75%  '

$ws.Range("E29").Value = '  +5.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.09'
$ws.Range("E30").Value = '  +2.34%  '

$ws.Range("E31").Value = '  +3.29%  '

$ws.Range("E32").Value = '  +2.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.57'
$ws.Range("E33").Value = '  +1.90%  '

$ws.Range("E34").Value = '  +15.77%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.11%  '

$ws.Range("E36").Value = '  -0.84%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.16'
$ws.Range("E37").Value = '  +6.46%  '

$ws.Range("E38").Value = '  +0.55%  '

$ws.Range("E39").Value = '  +5.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.31'
$ws.Range("E40").Value = '  +2.32%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '148.77'
$ws.Range("E41").Value = '  -0.59%  '

$ws.Range("E42").Value = '  +1.30%  '

$ws.Range("E43").Value = '  +2.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '285.30'
$ws.Range("E44").Value = '  +3.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0934'
$ws.Range("E45").Value = '  +1.55%  '

$ws.Range("E46").Value = '  +2.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.18'
$ws.Range("E47").Value = '  +7.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.559'
$ws.Range("E48").Value = '  +1.06%  '

$ws.Range("E49").Value = '  +2.97%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.56'
$ws.Range("E50").Value = '  +3.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.382'
$ws.Range("E51").Value = '  +1.50%  '
